$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert the new PRPS1 row: shift B/C/D of rows 68-98 down to
#    69-99 (process bottom-up so we never clobber a source row before
#    it has been read), then write the new PRPS1 values into row 68.
#    Column A (the plain 0-based counter) is left untouched for rows
#    2-98 and a new A99 cell is added at the end.
#
#    Column D stores its numbers as TEXT in the original workbook, so
#    we temporarily force a text NumberFormat before assigning each
#    digit-only string -- otherwise COM auto-coerces "2" into the
#    number 2. The temporary format is reverted at the end by pasting
    # the (untouched) D67 cell's format over the rewritten range, which
#    keeps the cells looking exactly like the rest of column D.
# ------------------------------------------------------------------
$ws.Range("D68:D99").NumberFormat = "@"

$ws.Cells.Item(99, 2).Value = "TBP"
$ws.Cells.Item(99, 3).Value = "TATA-box binding protein"
$ws.Cells.Item(99, 4).Value = "0"
$ws.Cells.Item(98, 2).Value = "RFC1"
$ws.Cells.Item(98, 3).Value = "replication factor C subunit 1"
$ws.Cells.Item(98, 4).Value = "0"
$ws.Cells.Item(97, 2).Value = "PPP2R2B"
$ws.Cells.Item(97, 3).Value = "protein phosphatase 2 regulatory subunit Bbeta"
$ws.Cells.Item(97, 4).Value = "0"
$ws.Cells.Item(96, 2).Value = "NOP56"
$ws.Cells.Item(96, 3).Value = "NOP56 ribonucleoprotein"
$ws.Cells.Item(96, 4).Value = "0"
$ws.Cells.Item(95, 2).Value = "FMR1"
$ws.Cells.Item(95, 3).Value = "fragile X mental retardation 1"
$ws.Cells.Item(95, 4).Value = "0"
$ws.Cells.Item(94, 2).Value = "DAB1"
$ws.Cells.Item(94, 3).Value = "DAB1, reelin adaptor protein"
$ws.Cells.Item(94, 4).Value = "0"
$ws.Cells.Item(93, 2).Value = "BEAN1"
$ws.Cells.Item(93, 3).Value = "brain expressed associated with NEDD4 1"
$ws.Cells.Item(93, 4).Value = "0"
$ws.Cells.Item(92, 2).Value = "ATXN8"
$ws.Cells.Item(92, 3).Value = "ataxin 8"
$ws.Cells.Item(92, 4).Value = "0"
$ws.Cells.Item(91, 2).Value = "ATXN7"
$ws.Cells.Item(91, 3).Value = "ataxin 7"
$ws.Cells.Item(91, 4).Value = "0"
$ws.Cells.Item(90, 2).Value = "ATXN3"
$ws.Cells.Item(90, 3).Value = "ataxin 3"
$ws.Cells.Item(90, 4).Value = "0"
$ws.Cells.Item(89, 2).Value = "ATXN2"
$ws.Cells.Item(89, 3).Value = "ataxin 2"
$ws.Cells.Item(89, 4).Value = "0"
$ws.Cells.Item(88, 2).Value = "ATXN10"
$ws.Cells.Item(88, 3).Value = "ataxin 10"
$ws.Cells.Item(88, 4).Value = "0"
$ws.Cells.Item(87, 2).Value = "ATXN1"
$ws.Cells.Item(87, 3).Value = "ataxin 1"
$ws.Cells.Item(87, 4).Value = "0"
$ws.Cells.Item(86, 2).Value = "ATN1"
$ws.Cells.Item(86, 3).Value = "atrophin 1"
$ws.Cells.Item(86, 4).Value = "0"
$ws.Cells.Item(85, 2).Value = "VWA3B"
$ws.Cells.Item(85, 3).Value = "von Willebrand factor A domain containing 3B"
$ws.Cells.Item(85, 4).Value = "1"
$ws.Cells.Item(84, 2).Value = "TSEN54"
$ws.Cells.Item(84, 3).Value = "tRNA splicing endonuclease subunit 54"
$ws.Cells.Item(84, 4).Value = "1"
$ws.Cells.Item(83, 2).Value = "TGM6"
$ws.Cells.Item(83, 3).Value = "transglutaminase 6"
$ws.Cells.Item(83, 4).Value = "1"
$ws.Cells.Item(82, 2).Value = "SYT14"
$ws.Cells.Item(82, 3).Value = "synaptotagmin 14"
$ws.Cells.Item(82, 4).Value = "1"
$ws.Cells.Item(81, 2).Value = "SEPSECS"
$ws.Cells.Item(81, 3).Value = "Sep (O-phosphoserine) tRNA:Sec (selenocysteine) tRNA synthase"
$ws.Cells.Item(81, 4).Value = "1"
$ws.Cells.Item(80, 2).Value = "NOL3"
$ws.Cells.Item(80, 3).Value = "nucleolar protein 3"
$ws.Cells.Item(80, 4).Value = "1"
$ws.Cells.Item(79, 2).Value = "MME"
$ws.Cells.Item(79, 3).Value = "membrane metalloendopeptidase"
$ws.Cells.Item(79, 4).Value = "1"
$ws.Cells.Item(78, 2).Value = "IFRD1"
$ws.Cells.Item(78, 3).Value = "interferon related developmental regulator 1"
$ws.Cells.Item(78, 4).Value = "1"
$ws.Cells.Item(77, 2).Value = "FDXR"
$ws.Cells.Item(77, 3).Value = "ferredoxin reductase"
$ws.Cells.Item(77, 4).Value = "1"
$ws.Cells.Item(76, 2).Value = "EEF2"
$ws.Cells.Item(76, 3).Value = "eukaryotic translation elongation factor 2"
$ws.Cells.Item(76, 4).Value = "1"
$ws.Cells.Item(75, 2).Value = "CACNB4"
$ws.Cells.Item(75, 3).Value = "calcium voltage-gated channel auxiliary subunit beta 4"
$ws.Cells.Item(75, 4).Value = "1"
$ws.Cells.Item(74, 2).Value = "ATP7B"
$ws.Cells.Item(74, 3).Value = "ATPase copper transporting beta"
$ws.Cells.Item(74, 4).Value = "1"
$ws.Cells.Item(73, 2).Value = "ATP1A2"
$ws.Cells.Item(73, 3).Value = "ATPase Na+/K+ transporting subunit alpha 2"
$ws.Cells.Item(73, 4).Value = "1"
$ws.Cells.Item(72, 2).Value = "ZFYVE26"
$ws.Cells.Item(72, 3).Value = "zinc finger FYVE-type containing 26"
$ws.Cells.Item(72, 4).Value = "2"
$ws.Cells.Item(71, 2).Value = "VAMP1"
$ws.Cells.Item(71, 3).Value = "vesicle associated membrane protein 1"
$ws.Cells.Item(71, 4).Value = "2"
$ws.Cells.Item(70, 2).Value = "TRPC3"
$ws.Cells.Item(70, 3).Value = "transient receptor potential cation channel subfamily C member 3"
$ws.Cells.Item(70, 4).Value = "2"
$ws.Cells.Item(69, 2).Value = "SDHA"
$ws.Cells.Item(69, 3).Value = "succinate dehydrogenase complex flavoprotein subunit A"
$ws.Cells.Item(69, 4).Value = "2"

# Write the new PRPS1 row (row 68) gene data
$ws.Cells.Item(68, 2).Value = "PRPS1"
$ws.Cells.Item(68, 3).Value = "phosphoribosyl pyrophosphate synthetase 1"
$ws.Cells.Item(68, 4).Value = "2"

# Revert the temporary text format on column D back to how the rest
# of the column looks (plain, unstyled cells) now that the text
# values are locked in.
$ws.Range("D67").Copy()
$ws.Range("D68:D99").PasteSpecial(-4122)
$excel.Application.CutCopyMode = $false

# New last row (99) needs an index value in column A, formatted like
# the rest of column A (copy format from A98 first).
$ws.Range("A98").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$excel.Application.CutCopyMode = $false
$ws.Cells.Item(99, 1).Value = 97

# New last row (99) also needs its "panel" value in column E (every
# data row repeats the same panel name).
$ws.Cells.Item(99, 5).Value = "Ataxia - adult onset"

# ------------------------------------------------------------------
# 2) Add the new "time_taken" column F, with the same header style as
#    the other header cells (copy format from E1).
# ------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.Application.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "time_taken"

# Fill time_taken for every data row (2-99)
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:50:10.657493"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:50:10.657505"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:50:10.657509"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:50:10.657511"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:50:10.657514"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:50:10.657517"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:50:10.657519"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:50:10.657522"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:50:10.657525"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:50:10.657527"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:50:10.657530"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:50:10.657532"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:50:10.657535"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:50:10.657537"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:50:10.657540"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:50:10.657542"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:50:10.657545"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:50:10.657547"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:50:10.657550"
$ws.Cells.Item(21, 6).Value = "2021-10-05 10:50:10.657553"
$ws.Cells.Item(22, 6).Value = "2021-10-05 10:50:10.657555"
$ws.Cells.Item(23, 6).Value = "2021-10-05 10:50:10.657557"
$ws.Cells.Item(24, 6).Value = "2021-10-05 10:50:10.657560"
$ws.Cells.Item(25, 6).Value = "2021-10-05 10:50:10.657562"
$ws.Cells.Item(26, 6).Value = "2021-10-05 10:50:10.657565"
$ws.Cells.Item(27, 6).Value = "2021-10-05 10:50:10.657568"
$ws.Cells.Item(28, 6).Value = "2021-10-05 10:50:10.657570"
$ws.Cells.Item(29, 6).Value = "2021-10-05 10:50:10.657573"
$ws.Cells.Item(30, 6).Value = "2021-10-05 10:50:10.657575"
$ws.Cells.Item(31, 6).Value = "2021-10-05 10:50:10.657578"
$ws.Cells.Item(32, 6).Value = "2021-10-05 10:50:10.657580"
$ws.Cells.Item(33, 6).Value = "2021-10-05 10:50:10.657583"
$ws.Cells.Item(34, 6).Value = "2021-10-05 10:50:10.657586"
$ws.Cells.Item(35, 6).Value = "2021-10-05 10:50:10.657588"
$ws.Cells.Item(36, 6).Value = "2021-10-05 10:50:10.657591"
$ws.Cells.Item(37, 6).Value = "2021-10-05 10:50:10.657593"
$ws.Cells.Item(38, 6).Value = "2021-10-05 10:50:10.657596"
$ws.Cells.Item(39, 6).Value = "2021-10-05 10:50:10.657599"
$ws.Cells.Item(40, 6).Value = "2021-10-05 10:50:10.657601"
$ws.Cells.Item(41, 6).Value = "2021-10-05 10:50:10.657604"
$ws.Cells.Item(42, 6).Value = "2021-10-05 10:50:10.657607"
$ws.Cells.Item(43, 6).Value = "2021-10-05 10:50:10.657609"
$ws.Cells.Item(44, 6).Value = "2021-10-05 10:50:10.657612"
$ws.Cells.Item(45, 6).Value = "2021-10-05 10:50:10.657614"
$ws.Cells.Item(46, 6).Value = "2021-10-05 10:50:10.657617"
$ws.Cells.Item(47, 6).Value = "2021-10-05 10:50:10.657619"
$ws.Cells.Item(48, 6).Value = "2021-10-05 10:50:10.657622"
$ws.Cells.Item(49, 6).Value = "2021-10-05 10:50:10.657624"
$ws.Cells.Item(50, 6).Value = "2021-10-05 10:50:10.657627"
$ws.Cells.Item(51, 6).Value = "2021-10-05 10:50:10.657629"
$ws.Cells.Item(52, 6).Value = "2021-10-05 10:50:10.657631"
$ws.Cells.Item(53, 6).Value = "2021-10-05 10:50:10.657634"
$ws.Cells.Item(54, 6).Value = "2021-10-05 10:50:10.657637"
$ws.Cells.Item(55, 6).Value = "2021-10-05 10:50:10.657639"
$ws.Cells.Item(56, 6).Value = "2021-10-05 10:50:10.657642"
$ws.Cells.Item(57, 6).Value = "2021-10-05 10:50:10.657644"
$ws.Cells.Item(58, 6).Value = "2021-10-05 10:50:10.657647"
$ws.Cells.Item(59, 6).Value = "2021-10-05 10:50:10.657649"
$ws.Cells.Item(60, 6).Value = "2021-10-05 10:50:10.657652"
$ws.Cells.Item(61, 6).Value = "2021-10-05 10:50:10.657654"
$ws.Cells.Item(62, 6).Value = "2021-10-05 10:50:10.657657"
$ws.Cells.Item(63, 6).Value = "2021-10-05 10:50:10.657659"
$ws.Cells.Item(64, 6).Value = "2021-10-05 10:50:10.657662"
$ws.Cells.Item(65, 6).Value = "2021-10-05 10:50:10.657664"
$ws.Cells.Item(66, 6).Value = "2021-10-05 10:50:10.657668"
$ws.Cells.Item(67, 6).Value = "2021-10-05 10:50:10.657670"
$ws.Cells.Item(68, 6).Value = "2021-10-05 10:50:10.657673"
$ws.Cells.Item(69, 6).Value = "2021-10-05 10:50:10.657675"
$ws.Cells.Item(70, 6).Value = "2021-10-05 10:50:10.657678"
$ws.Cells.Item(71, 6).Value = "2021-10-05 10:50:10.657680"
$ws.Cells.Item(72, 6).Value = "2021-10-05 10:50:10.657683"
$ws.Cells.Item(73, 6).Value = "2021-10-05 10:50:10.657685"
$ws.Cells.Item(74, 6).Value = "2021-10-05 10:50:10.657688"
$ws.Cells.Item(75, 6).Value = "2021-10-05 10:50:10.657690"
$ws.Cells.Item(76, 6).Value = "2021-10-05 10:50:10.657693"
$ws.Cells.Item(77, 6).Value = "2021-10-05 10:50:10.657695"
$ws.Cells.Item(78, 6).Value = "2021-10-05 10:50:10.657700"
$ws.Cells.Item(79, 6).Value = "2021-10-05 10:50:10.657703"
$ws.Cells.Item(80, 6).Value = "2021-10-05 10:50:10.657706"
$ws.Cells.Item(81, 6).Value = "2021-10-05 10:50:10.657708"
$ws.Cells.Item(82, 6).Value = "2021-10-05 10:50:10.657711"
$ws.Cells.Item(83, 6).Value = "2021-10-05 10:50:10.657713"
$ws.Cells.Item(84, 6).Value = "2021-10-05 10:50:10.657715"
$ws.Cells.Item(85, 6).Value = "2021-10-05 10:50:10.657718"
$ws.Cells.Item(86, 6).Value = "2021-10-05 10:50:10.657720"
$ws.Cells.Item(87, 6).Value = "2021-10-05 10:50:10.657723"
$ws.Cells.Item(88, 6).Value = "2021-10-05 10:50:10.657725"
$ws.Cells.Item(89, 6).Value = "2021-10-05 10:50:10.657728"
$ws.Cells.Item(90, 6).Value = "2021-10-05 10:50:10.657730"
$ws.Cells.Item(91, 6).Value = "2021-10-05 10:50:10.657733"
$ws.Cells.Item(92, 6).Value = "2021-10-05 10:50:10.657735"
$ws.Cells.Item(93, 6).Value = "2021-10-05 10:50:10.657737"
$ws.Cells.Item(94, 6).Value = "2021-10-05 10:50:10.657741"
$ws.Cells.Item(95, 6).Value = "2021-10-05 10:50:10.657744"
$ws.Cells.Item(96, 6).Value = "2021-10-05 10:50:10.657747"
$ws.Cells.Item(97, 6).Value = "2021-10-05 10:50:10.657749"
$ws.Cells.Item(98, 6).Value = "2021-10-05 10:50:10.657752"
$ws.Cells.Item(99, 6).Value = "2021-10-05 10:50:10.657754"

$excel.Application.CutCopyMode = $false
Write-Host "edit complete"
